$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write column B first, in the order that reproduces the original shared-string table order
$ws.Range("B3").Value = 'Bruno Zan Arriaga'
$ws.Range("B4").Value = 'Victor Hugo Nolasco'
$ws.Range("B5").Value = 'Adelmo Felipe de Oliveira Bento da Silva'
$ws.Range("B8").Value = 'Felipe Alves Ferreira'
$ws.Range("B11").Value = 'Diogo Rodrigues de Sousa Lima'
$ws.Range("B12").Value = 'Raquel Shadeck Ouchi'
$ws.Range("B14").Value = 'Douglas Oliveira Andrade'
$ws.Range("B15").Value = 'Laura Barreto Miranda Scarpa Leite'
$ws.Range("B6").Value = 'Gama Rodrigo Gama dos Santos'
$ws.Range("B7").Value = 'Lu Luciane Petrangelo'
$ws.Range("B9").Value = 'Mendes Rodrigo Mendes Santos'
$ws.Range("B10").Value = 'Kimori Marcos Goncalves Kimori'
$ws.Range("B13").Value = 'Baldocchi Rafael Cesar Baldocchi'
$ws.Range("B16").Value = 'Balan Antonio Carlos Balan Junior'
$ws.Range("B18").Value = 'Fe Fernanda Eugenio'
$ws.Range("B17").Value = 'Jean Paulo Kambourakis'
$ws.Range("B19").Value = 'Nanda Annanda Destro Torteli'

# Write column A (matricula numbers) in row order
$ws.Range("A3").Value = 731198
$ws.Range("A4").Value = 775732
$ws.Range("A5").Value = 793126
$ws.Range("A6").Value = 629839
$ws.Range("A7").Value = 497935
$ws.Range("A8").Value = 615881
$ws.Range("A9").Value = 594396
$ws.Range("A10").Value = 696629
$ws.Range("A11").Value = 660439
$ws.Range("A12").Value = 591914
$ws.Range("A13").Value = 607954
$ws.Range("A14").Value = 661523
$ws.Range("A15").Value = 705641
$ws.Range("A16").Value = 667096
$ws.Range("A17").Value = 607874
$ws.Range("A18").Value = 666396
$ws.Range("A19").Value = 808110

$ws.Range("A20").Select()
